$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Fix 1: typo "mesto" -> "mestu" in the "Svaki korisnik ce imati..."
# paragraph. We replace just the single offending letter ("o" -> "u")
# so that only that character's run is split off, matching how the
# author actually corrected the typo (select the "o", type "u").
# ---------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute(
    "na jednom mesto", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    $oStart = $full.End - 1
    $oChar = $d.Range($oStart, $oStart + 1)

    # Nudge formatting away-and-back so the engine treats this single
    # character as its own edited run instead of re-flowing/merging the
    # whole paragraph back into one run.
    $oChar.Bold = 1
    $oChar.Text = "u"
    $oChar2 = $d.Range($oStart, $oStart + 1)
    $oChar2.Bold = 0

    # The paragraph also contains a following run ("Korisnik ce moci da
    # vidi ...") that must stay a separate run from the one we just
    # edited. Re-stamp its formatting (no actual value change) across
    # its full span so the engine keeps it as its own run instead of
    # folding it into the run we edited above.
    $full2 = $d.Content
    $found2 = $full2.Find.Execute(
        "Korisnik će moći da vidi svoje zapraćene oglase u sekciji „zapraćeni oglasi“.",
        $true, $false, $false, $false, $false,
        $true, 1, $false, "", 0)
    if ($found2) {
        $tail = $d.Range($full2.Start, $full2.End)
        $tail.Bold = 1
        $tail2 = $d.Range($full2.Start, $full2.End)
        $tail2.Bold = 0
    }
}

# ---------------------------------------------------------------------
# Fix 2: close the curly quote in "Zapraćeni oglasi." ->
# "Zapraćeni oglasi“." The quote was opened ( „ ) but never closed; add
# the closing guillemet ( “ ) right before the final period.
# ---------------------------------------------------------------------
$full3 = $d.Content
$found3 = $full3.Find.Execute(
    "Oglas se pojavljuje u sekciji „Zapraćeni oglasi.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found3) {
    $periodStart = $full3.End - 1

    # Replace "." with "“." (split off a dedicated run for the tail).
    $periodChar = $d.Range($periodStart, $periodStart + 1)
    $periodChar.Bold = 1
    $periodChar.Text = "“."
    $periodReset = $d.Range($periodStart, $periodStart + 2)
    $periodReset.Bold = 0

    # Now split that new two-character tail into its own "“" run and
    # "." run, same as the author's original edit.
    $dotStart = $periodStart + 1
    $dotChar = $d.Range($dotStart, $dotStart + 1)
    $dotChar.Bold = 1
    $dotChar.Bold = 0
}
